# Auto-generated edit script: updates cryptos list values (price/volume) per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.156.49"
$ws.Range("E2").Value = "  +7.83%  "
$ws.Range("D3").Value = "1.588.24"
$ws.Range("E3").Value = "  +7.87%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'0.9900"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.27%  "
$ws.Range("D6").Value = "'297.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.41%  "
$ws.Range("D7").Value = "'0.3625"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").Value = "'0.3341"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +8.67%  "
$ws.Range("D9").Value = "'41.30"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.15%  "
$ws.Range("D10").Value = "'1.117"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.75%  "
$ws.Range("D11").Value = "'0.06942"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.28%  "
$ws.Range("D12").Value = "'1.003"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D13").Value = "'19.38"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.49%  "
$ws.Range("D14").Value = "'5.813"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.24%  "
$ws.Range("D15").Value = "'6.524"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.70%  "
$ws.Range("D16").Value = "'0.9917"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.45%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "1.593.00"
$ws.Range("E17").Value = "  +8.14%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.00001062"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.32%  "
$ws.Range("D19").Value = "'0.06555"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +10.44%  "
$ws.Range("D20").Value = "'75.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +10.18%  "
$ws.Range("D21").Value = "'15.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +8.57%  "
$ws.Range("D22").Value = "'5.910"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.45%  "
$ws.Range("D23").Value = "'11.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.93%  "
$ws.Range("D24").Value = "22.175.83"
$ws.Range("E24").Value = "  +7.86%  "
$ws.Range("D25").Value = "'2.376"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.88%  "
$ws.Range("D26").Value = "'2.490"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +16.80%  "
$ws.Range("D27").Value = "'148.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Value = "'19.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +11.45%  "
$ws.Range("D29").Value = "1.761.45"
$ws.Range("E29").Value = "  +7.88%  "
$ws.Range("D30").Value = "'121.43"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.51%  "
$ws.Range("D31").Value = "'3.973"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.83%  "
$ws.Range("D32").Value = "'5.883"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +18.59%  "
$ws.Range("D33").Value = "'0.9164"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +13.66%  "
$ws.Range("D34").Value = "'0.08159"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.71%  "
$ws.Range("D35").Value = "'1.614"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.52%  "
$ws.Range("D36").Value = "'11.65"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +11.62%  "
$ws.Range("D37").Value = "'5.102"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.52%  "
$ws.Range("D38").Value = "'1.234"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.67%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.06001"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.12%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'8.316"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +11.83%  "
$ws.Range("D41").Value = "'0.02173"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.75%  "
$ws.Range("D42").Value = "'0.1978"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.63%  "
$ws.Range("D43").Value = "'0.9917"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.37%  "
$ws.Range("D44").Value = "'0.5764"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.95%  "
$ws.Range("D45").Value = "'3.758"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.82%  "
$ws.Range("D46").Value = "'12.94"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.38%  "
$ws.Range("D47").Value = "'124.93"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.31%  "
$ws.Range("D48").Value = "'0.5551"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.57%  "
$ws.Range("D49").Value = "'1.935"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.23%  "
$ws.Range("D50").Value = "'0.06703"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.51%  "
$ws.Range("D51").Value = "'72.38"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +8.06%  "
